# Applies the "add partner name to the order form" edit to the Piql order
# form workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block -----------------------------------------------------
# Date: field gets re-stamped (time-of-day moves forward).
$ws.Range("G4").Value = 44079.5792554584

# Piql Partner: name is now filled in.
$ws.Range("G5").Value = "alfrtruj"

# Customer name changes.
$ws.Range("G7").Value = "El enano"

# Address/comments field changes.
$ws.Range("F10").Value = "el perro"

# --- Order lines --------------------------------------------------------
# 1. Offline Storage
#    piqlConnect (only piqlFilm) line - cleared out.
$ws.Range("F18").Value = $null
$ws.Range("G18").Value = $null
$ws.Range("H18").Value = $null

#    Digital (GB) line - now has quantity/price/total.
$ws.Range("F19").Value = 450
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 4950

#    Visual (pages) line - cleared out (including Pages/frame helper cell).
$ws.Range("E20").Value = $null
$ws.Range("F20").Value = $null
$ws.Range("G20").Value = $null
$ws.Range("H20").Value = $null

# 2. Online Storage (GB) - piqlConnect
#    Section row now carries the quantity/price/total.
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 8880
$ws.Range("H21").Value = 8880

#    Online Storage (GB) line now carries payment term + totals.
$ws.Range("E22").Value = "yearly"
$ws.Range("F22").Value = 7890
$ws.Range("G22").Value = 0.576
$ws.Range("H22").Value = 3968.64

# 6. Shipment cost
#    Reels count doubles, total follows.
$ws.Range("E32").Value = 4
$ws.Range("H32").Value = 120

# TOTAL row.
$ws.Range("H33").Value = 17798.64

# Total to pay from the second term.
$ws.Range("H34").Value = 12848.64
